# HvilleHappyHours.xlsx — "Add files via upload / Updated Happy Hours and Specials"
#
# Adds new Happy-Hour / Specials text to the Shenanigans row (12), and to the
# Jolly Fox Club (18) / Time Out Karaoke (19) rows, then nudges a few
# formatting bits (row heights, one column width, the saved selection) to
# match the refreshed sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shenanigans (row 12): fill in the per-day specials -------------------
$ws.Range("D12").Value = "ALL NIGHT                                $0.25 Wells and Draft"
$ws.Range("E12").Value = "ALL NIGHT                                          $1.50 Wells/Domestic Longnecks                                               $250 You-Call-Its"
$ws.Range("F12").Value = "ALL NIGHT                          $1.75 Domestic Longnecks                                                     8PM-11PM                               $0.50 Wells/Drafts"
$ws.Range("G12").Value = "8PM-11PM                             $1.00 Wells/Dom Longnecks                                           $3.00 You-Call-Its"
$ws.Range("H12").Value = "ALL NIGHT                                 $1.00 Wells/Drafts                           $1.50 Dom longnecks"

# --- Jolly Fox Club (row 18): Wednesday special ----------------------------
$ws.Range("D18").Value = "ALL NIGHT                                $0.25 Wells and Draft"

# --- Time Out Karaoke (row 19): Thursday special + notes column -----------
$ws.Range("E19").Value = "ALL NIGHT                                First Drink reg price,              2nd Drink `$3 OFF                          (as low as `$1)                                          `$3.75 Shots"
$ws.Range("I19").Value = "**Special Vary on Events**"

# --- Jolly Fox Club (row 18) notes column ----------------------------------
$ws.Range("I18").Value = "**Specials Vary Daily**                                                               **Check Twitter for Updates**"

# --- Row heights to fit the newly-added wrapped text -----------------------
$ws.Rows(12).RowHeight = 53.25
$ws.Rows(18).RowHeight = 30

# --- Column D widened slightly to match the refreshed layout --------------
$ws.Columns(4).ColumnWidth = 25.85

# --- Leave the saved selection where the edits finished --------------------
$ws.Range("I18").Select()
